# Weekly update: insert a new "Haba" price record for
# Comercializadora del Agro de Limarí (row 95), shifting the
# existing historical rows (95-116) down by one (96-117).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 95; this pushes old rows 95..116 to 96..117
$ws.Rows.Item(95).Insert()

# Populate the new row 95 with the latest weekly record
$ws.Cells.Item(95, 1).Value = 2
$ws.Cells.Item(95, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(95, 3).Value = "Coquimbo"
$ws.Cells.Item(95, 4).Value = 45211
$ws.Cells.Item(95, 5).Value = 4
$ws.Cells.Item(95, 6).Value = 100112026
$ws.Cells.Item(95, 7).Value = "Haba"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 1000
$ws.Cells.Item(95, 11).Value = 7500
$ws.Cells.Item(95, 12).Value = 8000
$ws.Cells.Item(95, 13).Value = 7750
$ws.Cells.Item(95, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(95, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(95, 16).Value = 310
$ws.Cells.Item(95, 17).Value = 25
$ws.Cells.Item(95, 18).Value = "Hortaliza"
